# Add a new "walk in closet" feature row to Sheet1 and update the
# active cell selection, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: index number + new feature name ("walk in closet" is appended
# as a brand new shared string).
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "walk in closet"

# Move/record the active selection as in the diff (B5 -> F12).
$ws.Range("F12").Select()
